$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1559.375
$ws.Range("J17").Value = 1559.375
$ws.Range("L17").Value = 4678.125
$ws.Range("N17").Value = -5014.125
$ws.Range("H28").Value = 1369.35
$ws.Range("J28").Value = 4339
$ws.Range("L28").Value = 4339
$ws.Range("N28").Value = -5309
$ws.Range("H39").Value = 2750.0667
$ws.Range("I39").Value = 1200.3
$ws.Range("J39").Value = 5849.6
$ws.Range("K39").Value = 3600.9
$ws.Range("L39").Value = 17548.8
$ws.Range("M39").Value = -3304.9
$ws.Range("N39").Value = -18140.8
$ws.Range("H61").Value = 1915
$ws.Range("I61").Value = 1915
$ws.Range("K61").Value = 5745
$ws.Range("M61").Value = -5573
$ws.Range("H69").Value = 19624.125
$ws.Range("I69").Value = 18496.5
$ws.Range("K69").Value = 55489.5
$ws.Range("M69").Value = -54615.5
$ws.Range("H72").Value = 19624.125
$ws.Range("I72").Value = 18496.5
$ws.Range("K72").Value = 166468.5
$ws.Range("M72").Value = -162100.5
$ws.Range("H82").Value = 1024.5
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1024.5
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 5040.6665
$ws.Range("J86").Value = 7332
$ws.Range("L86").Value = 7332
$ws.Range("N86").Value = -9578
$ws.Range("H89").Value = 5040.6665
$ws.Range("J89").Value = 7332
$ws.Range("L89").Value = 36660
$ws.Range("N89").Value = -47892
$ws.Range("H100").Value = 7201
$ws.Range("I100").Value = 6268
$ws.Range("K100").Value = 6268
$ws.Range("M100").Value = -5727
$ws.Range("H106").Value = 8419.666999999999
$ws.Range("I106").Value = 8419.666999999999
$ws.Range("K106").Value = 8419.666999999999
$ws.Range("M106").Value = -7788.666999999999
$ws.Range("H129").Value = 1904.2858
$ws.Range("I129").Value = 1904.2858
$ws.Range("K129").Value = 5712.857400000001
$ws.Range("M129").Value = -712.8574000000008
$ws.Range("H138").Value = 2647.8572
$ws.Range("I138").Value = 1219.5385
$ws.Range("J138").Value = 3079.6743
$ws.Range("K138").Value = 3658.6155
$ws.Range("L138").Value = 9239.0229
$ws.Range("M138").Value = 1481.3845
$ws.Range("N138").Value = -19519.0229

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3065.389
$ws.Range("I2").Value = 3057.25
$ws.Range("K2").Value = 3057.25
$ws.Range("M2").Value = -2944.25
$ws.Range("H45").Value = 10825.88
$ws.Range("I45").Value = 13173.685
$ws.Range("K45").Value = 13173.685
$ws.Range("M45").Value = -12796.685
$ws.Range("H63").Value = 2075.9443
$ws.Range("I63").Value = 2048.5715
$ws.Range("K63").Value = 2048.5715
$ws.Range("M63").Value = -1362.5715
$ws.Range("H66").Value = 2075.9443
$ws.Range("I66").Value = 2048.5715
$ws.Range("K66").Value = 10242.8575
$ws.Range("M66").Value = -6810.8575
$ws.Range("H97").Value = 985.64703
$ws.Range("I97").Value = 709.9286
$ws.Range("K97").Value = 709.9286
$ws.Range("M97").Value = -213.9286
$ws.Range("H116").Value = 3065.389
$ws.Range("I116").Value = 3057.25
$ws.Range("K116").Value = 3057.25
$ws.Range("M116").Value = -763.25
$ws.Range("H128").Value = 9999.5
$ws.Range("J128").Value = 9999.5
$ws.Range("L128").Value = 9999.5
$ws.Range("N128").Value = -19959.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3065.389
$ws.Range("I3").Value = 3057.25
$ws.Range("K3").Value = 3057.25
$ws.Range("M3").Value = -2943.25
$ws.Range("H107").Value = 3422.6667
$ws.Range("I107").Value = 3300
$ws.Range("K107").Value = 3300
$ws.Range("M107").Value = -1380
$ws.Range("H134").Value = 2611.9355
$ws.Range("I134").Value = 2632.3333
$ws.Range("K134").Value = 7896.999899999999
$ws.Range("M134").Value = -5361.999899999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 9505.385
$ws.Range("J99").Value = 12414.333
$ws.Range("L99").Value = 12414.333
$ws.Range("N99").Value = -15410.333
$ws.Range("H126").Value = 9505.385
$ws.Range("J126").Value = 12414.333
$ws.Range("L126").Value = 37242.999
$ws.Range("N126").Value = -42182.999
$ws.Range("H133").Value = 73215.336
$ws.Range("J133").Value = 74799.2
$ws.Range("L133").Value = 74799.2
$ws.Range("N133").Value = -79859.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 4000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 4000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 12000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -13118
$ws.Range("H70").Value = 18102.4
$ws.Range("I70").Value = 13256
$ws.Range("J70").Value = 21333.334
$ws.Range("K70").Value = 39768
$ws.Range("L70").Value = 64000.00199999999
$ws.Range("M70").Value = -39453
$ws.Range("N70").Value = -64630.00199999999
$ws.Range("H73").Value = 18102.4
$ws.Range("I73").Value = 13256
$ws.Range("J73").Value = 21333.334
$ws.Range("K73").Value = 39768
$ws.Range("L73").Value = 64000.00199999999
$ws.Range("M73").Value = -38676
$ws.Range("N73").Value = -66184.00199999999
$ws.Range("H75").Value = 3876.625
$ws.Range("J75").Value = 3873.5715
$ws.Range("L75").Value = 11620.7145
$ws.Range("N75").Value = -13616.7145
$ws.Range("H78").Value = 3876.625
$ws.Range("J78").Value = 3873.5715
$ws.Range("L78").Value = 34862.1435
$ws.Range("N78").Value = -44846.1435
$ws.Range("H88").Value = 19999.666
$ws.Range("H91").Value = 19999.666
$ws.Range("H139").Value = 4050.2778
$ws.Range("I139").Value = 2350.7144
$ws.Range("K139").Value = 7052.1432
$ws.Range("M139").Value = -1912.1432
$ws.Range("H140").Value = 4668.9287
$ws.Range("I140").Value = 3761.2727
$ws.Range("K140").Value = 11283.8181
$ws.Range("M140").Value = -6103.8181

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5946.3335
$ws.Range("I70").Value = 5949.5
$ws.Range("K70").Value = 5949.5
$ws.Range("M70").Value = -5679.5
$ws.Range("H73").Value = 5946.3335
$ws.Range("I73").Value = 5949.5
$ws.Range("K73").Value = 5949.5
$ws.Range("M73").Value = -5013.5
$ws.Range("H97").Value = 355.0625
$ws.Range("I97").Value = 312.14285
$ws.Range("K97").Value = 312.14285
$ws.Range("M97").Value = 183.85715
$ws.Range("H122").Value = 2633.8572
$ws.Range("I122").Value = 1984.25
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 5952.75
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -3502.75
$ws.Range("N122").Value = -15400
$ws.Range("H126").Value = 2144.8
$ws.Range("I126").Value = 1350
$ws.Range("K126").Value = 4050
$ws.Range("M126").Value = -1580

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1820.5454
$ws.Range("I22").Value = 1896.5714
$ws.Range("J22").Value = 1687.5
$ws.Range("K22").Value = 1896.5714
$ws.Range("L22").Value = 1687.5
$ws.Range("M22").Value = -1601.5714
$ws.Range("N22").Value = -2277.5
$ws.Range("H27").Value = 1820.5454
$ws.Range("I27").Value = 1896.5714
$ws.Range("J27").Value = 1687.5
$ws.Range("K27").Value = 1896.5714
$ws.Range("L27").Value = 1687.5
$ws.Range("M27").Value = -1789.5714
$ws.Range("N27").Value = -1901.5
$ws.Range("H122").Value = 12237.333
$ws.Range("J122").Value = 2980
$ws.Range("L122").Value = 8940
$ws.Range("N122").Value = -13840
$ws.Range("H125").Value = 89998.836
$ws.Range("J125").Value = 89998.836
$ws.Range("L125").Value = 89998.836
$ws.Range("N125").Value = -99838.836

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4056.7334
$ws.Range("I107").Value = 3596.2307
$ws.Range("K107").Value = 10788.6921
$ws.Range("M107").Value = -8868.6921
$ws.Range("H122").Value = 2538.7273
$ws.Range("I122").Value = 2325.6667
$ws.Range("J122").Value = 3497.5
$ws.Range("K122").Value = 6977.000100000001
$ws.Range("L122").Value = 10492.5
$ws.Range("M122").Value = -4527.000100000001
$ws.Range("N122").Value = -15392.5
$ws.Range("H126").Value = 2166.25
$ws.Range("I126").Value = 2166.25
$ws.Range("K126").Value = 6498.75
$ws.Range("M126").Value = -4028.75
$ws.Range("H131").Value = 109678.5
$ws.Range("J131").Value = 109678.5
$ws.Range("L131").Value = 109678.5
$ws.Range("N131").Value = -119758.5
$ws.Range("H132").Value = 4013.9395
$ws.Range("I132").Value = 3312.9644
$ws.Range("K132").Value = 9938.893199999999
$ws.Range("M132").Value = -7408.893199999999
$ws.Range("H136").Value = 1317.45
$ws.Range("I136").Value = 1135.9744
$ws.Range("K136").Value = 3407.9232
$ws.Range("M136").Value = -857.9232000000002
